$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename "2GB Ram server" -> "2GB Ram server vm" ---
$ws.Range("A1").Value = "2GB Ram server vm"

# --- Rows 14-20: col B becomes full date-time (day 42804) with m/d/yy h:mm format ---
$ws.Range("B14").Value = 42804.071527777778
$ws.Range("B14").NumberFormat = "m/d/yy h:mm"
$ws.Range("B15").Value = 42804.071527777778
$ws.Range("B15").NumberFormat = "m/d/yy h:mm"
$ws.Range("B16").Value = 42804.071527777778
$ws.Range("B16").NumberFormat = "m/d/yy h:mm"
$ws.Range("B17").Value = 42804.071527777778
$ws.Range("B17").NumberFormat = "m/d/yy h:mm"
$ws.Range("B18").Value = 42804.071527777778
$ws.Range("B18").NumberFormat = "m/d/yy h:mm"
$ws.Range("B19").Value = 42804.071527777778
$ws.Range("B19").NumberFormat = "m/d/yy h:mm"
$ws.Range("B20").Value = 42804.071527777778
$ws.Range("B20").NumberFormat = "m/d/yy h:mm"

# --- New column C width ---
$ws.Columns.Item(3).ColumnWidth = 20.43

# --- New data blocks: rows 22-41 and 44-63 ("FEC 2GB 1.99GHz quad" server runs) ---
$ws.Range("A22").Value = "FEC 2GB 1.99GHz quad"
$ws.Range("A23").Value = "vm1"
$ws.Range("B23").Value = 42804.40625
$ws.Range("B23").NumberFormat = "m/d/yy h:mm"
$ws.Range("A24").Value = "vm2"
$ws.Range("B24").Value = 42804.418055555558
$ws.Range("B24").NumberFormat = "m/d/yy h:mm"
$ws.Range("A25").Value = "vm3"
$ws.Range("B25").Value = 42804.432638888888
$ws.Range("B25").NumberFormat = "m/d/yy h:mm"
$ws.Range("C25").Value = 42804.448611111111
$ws.Range("C25").NumberFormat = "m/d/yy h:mm"
$ws.Range("A26").Value = "vm4"
$ws.Range("B26").Value = 42804.463194444441
$ws.Range("B26").NumberFormat = "m/d/yy h:mm"
$ws.Range("C26").NumberFormat = "h:mm"
$ws.Range("A27").Value = "vm5"
$ws.Range("B27").Value = 42804.463194444441
$ws.Range("B27").NumberFormat = "m/d/yy h:mm"
$ws.Range("A28").Value = "vm6"
$ws.Range("B28").Value = 42804.478472222225
$ws.Range("B28").NumberFormat = "m/d/yy h:mm"
$ws.Range("C28").Value = 0.48402777777777778
$ws.Range("C28").NumberFormat = "h:mm"
$ws.Range("A29").Value = "vm7"
$ws.Range("B29").Value = 42804.478472222225
$ws.Range("B29").NumberFormat = "m/d/yy h:mm"
$ws.Range("C29").Value = 0.48402777777777778
$ws.Range("C29").NumberFormat = "h:mm"
$ws.Range("A30").Value = "vm8"
$ws.Range("B30").Value = 42804.478472222225
$ws.Range("B30").NumberFormat = "m/d/yy h:mm"
$ws.Range("C30").Value = 0.48402777777777778
$ws.Range("C30").NumberFormat = "h:mm"
$ws.Range("A31").Value = "vm9"
$ws.Range("A32").Value = "vm10"
$ws.Range("A33").Value = "vm11"
$ws.Range("A34").Value = "vm12"
$ws.Range("A35").Value = "vm13"
$ws.Range("A36").Value = "vm14"
$ws.Range("A37").Value = "vm15"
$ws.Range("A38").Value = "vm16"
$ws.Range("A39").Value = "vm17"
$ws.Range("A40").Value = "vm18"
$ws.Range("A41").Value = "vm19"
$ws.Range("A44").Value = "FEC 2GB 1.99GHz quad"
$ws.Range("A45").Value = "vm1"
$ws.Range("B45").Value = 42804.520138888889
$ws.Range("B45").NumberFormat = "m/d/yy h:mm"
$ws.Range("A46").Value = "vm2"
$ws.Range("B46").Value = 42804.520138888889
$ws.Range("B46").NumberFormat = "m/d/yy h:mm"
$ws.Range("A47").Value = "vm3"
$ws.Range("B47").Value = 42804.520138888889
$ws.Range("B47").NumberFormat = "m/d/yy h:mm"
$ws.Range("C47").NumberFormat = "m/d/yy h:mm"
$ws.Range("A48").Value = "vm4"
$ws.Range("B48").Value = 42804.520138888889
$ws.Range("B48").NumberFormat = "m/d/yy h:mm"
$ws.Range("C48").NumberFormat = "h:mm"
$ws.Range("A49").Value = "vm5"
$ws.Range("B49").Value = 42804.520138888889
$ws.Range("B49").NumberFormat = "m/d/yy h:mm"
$ws.Range("A50").Value = "vm6"
$ws.Range("B50").Value = 42804.529861111114
$ws.Range("B50").NumberFormat = "m/d/yy h:mm"
$ws.Range("C50").NumberFormat = "h:mm"
$ws.Range("A51").Value = "vm7"
$ws.Range("B51").Value = 42804.529861111114
$ws.Range("B51").NumberFormat = "m/d/yy h:mm"
$ws.Range("C51").NumberFormat = "h:mm"
$ws.Range("A52").Value = "vm8"
$ws.Range("B52").Value = 42804.529861111114
$ws.Range("B52").NumberFormat = "m/d/yy h:mm"
$ws.Range("C52").NumberFormat = "h:mm"
$ws.Range("A53").Value = "vm9"
$ws.Range("A54").Value = "vm10"
$ws.Range("A55").Value = "vm11"
$ws.Range("A56").Value = "vm12"
$ws.Range("A57").Value = "vm13"
$ws.Range("A58").Value = "vm14"
$ws.Range("A59").Value = "vm15"
$ws.Range("A60").Value = "vm16"
$ws.Range("A61").Value = "vm17"
$ws.Range("A62").Value = "vm18"
$ws.Range("A63").Value = "vm19"

# --- View state (best effort) ---
$ws.Range("C50").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 28

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
